$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BF column holds the game-date string; it was off by one day (dates were
# recorded as "M-D-YYYY-YY" season strings instead of the real game date).
# Rewrite each as an ISO "YYYY-MM-DD" string.
#
# Force text entry (NumberFormat "@") so Excel doesn't reinterpret the
# literal "2014-05-29" as a date serial, then reset the style back to
# "Normal" so the cell keeps its original (default) formatting/style.
$dataRange = $ws.Range("BF2:BF31")
$dataRange.NumberFormat = "@"

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 58)
    $cell.Value = "2014-05-29"
    $cell.Style = "Normal"
}
